# hys 2 day 2
# Adds a "day2" column (D) to the "37" and "42" sheets, mirroring the
# existing "day1" column (C) layout/styling, and fills in the new values.

$wb = $excel.ActiveWorkbook

$sheet37Values = @{
    2  = 0.28189998865127563
    3  = 0.05119999870657921
    4  = 0.30640000104904175
    5  = 0.3197999894618988
    6  = 0.16899999976158142
    7  = 0.1623000055551529
    8  = 0.09920000284910202
    9  = 0.04989999905228615
    10 = 0.16269999742507935
    11 = 0.14010000228881836
    12 = 0.053199999034404755
    13 = 0.10840000212192535
    14 = 0.04769999906420708
    15 = 0.3425999879837036
    16 = 0.28110000491142273
    17 = 0.2694000005722046
    18 = 0.3303000032901764
    19 = 0.11079999804496765
    20 = 0.1875
    21 = 0.053300000727176666
    22 = 0.423799991607666
    23 = 0.05119999870657921
    24 = 0.31299999356269836
    25 = 0.15690000355243683
}

$sheet42Values = @{
    2  = 0.050200000405311584
    3  = 0.050999999046325684
    4  = 0.05130000039935112
    5  = 0.05040000006556511
    6  = 0.050700001418590546
    7  = 0.04969999939203262
    8  = 0.05139999836683273
    9  = 0.05299999937415123
    10 = 0.05009999871253967
    11 = 0.048700001090765
    12 = 0.0478999987244606
    13 = 0.04859999939799309
    14 = 0.04479999840259552
    15 = 0.04500000178813934
    16 = 0.04800000041723251
    17 = 0.04540000110864639
    18 = 0.05139999836683273
    19 = 0.05009999871253967
    20 = 0.0494999997317791
    21 = 0.048500001430511475
    22 = 0.04540000110864639
    23 = 0.04540000110864639
    24 = 0.04650000110268593
    25 = 0.045899998396635056
}

function Apply-Day2Column {
    param(
        $ws,
        [hashtable]$values
    )

    # Header: mirror C1's formatting for the new D1 "day2" header cell.
    $ws.Range("C1").Copy()
    $ws.Range("D1").PasteSpecial(-4122)
    $ws.Range("D1").Value = "day2"

    foreach ($row in $values.Keys | Sort-Object) {
        $cCell = $ws.Cells.Item($row, 3)
        $dCell = $ws.Cells.Item($row, 4)
        $cCell.Copy()
        $dCell.PasteSpecial(-4122)
        $dCell.Value = $values[$row]
    }
}

$ws37 = $wb.Worksheets.Item("37")
Apply-Day2Column $ws37 $sheet37Values

$ws42 = $wb.Worksheets.Item("42")
Apply-Day2Column $ws42 $sheet42Values
